# Fill in the remaining test-case rows of the "cost of gas" lab table and
# wire up the out:cost column with the real formula =(miles/MPG)*price,
# replacing the placeholder text "insert Excel formula here" in E2 and
# completing rows 3-5 (which were previously blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: "Zero miles, price is 1"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 30
$ws.Range("D3").Value = 1

# Row 4: "Number of miles is greater than MPG"
$ws.Range("B4").Value = 50
$ws.Range("C4").Value = 10
$ws.Range("D4").Value = 5

# Row 5: "Number of miles is less than MPG"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 60
$ws.Range("D5").Value = 3

# Write the calculated out:cost formula across the whole table column
# (E2:E5) in one shot so it is recognised as the table's calculated
# column and filled down automatically.
$lo = $ws.ListObjects.Item(1)
$col = $lo.ListColumns.Item(5)
$col.DataBodyRange.Formula = "=(B2/C2)*D2"

# Move the active selection off the table, as left by the author.
$ws.Range("E6").Select()
